# Update "想去人数" (interest count) figures in the "展览" and "全部类型" sheets,
# matching the refreshed data output committed to gh-pages.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): rows 2,3,4,5,9,10 in column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 332
$wsExhibit.Range("F3").Value = 100
$wsExhibit.Range("F4").Value = 494
$wsExhibit.Range("F5").Value = 4902
$wsExhibit.Range("F9").Value = 744
$wsExhibit.Range("F10").Value = 231

# Sheet "全部类型" (sheet4): rows 2,3,4,5,9,11 in column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 332
$wsAll.Range("F3").Value = 100
$wsAll.Range("F4").Value = 494
$wsAll.Range("F5").Value = 4902
$wsAll.Range("F9").Value = 744
$wsAll.Range("F11").Value = 231
